# Update the "2006-2018" history sheet with a new 2019 row at the top of
# the data (row 2), pushing all existing years down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2006-2018")

# Insert a new row above the current row 2 (year 2018), shifting the
# existing data (2018..2006) down to rows 3..15.
$ws.Rows("2:2").Insert()

# Fill in the new 2019 row.
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = "Andrés Gomez"
$ws.Range("C2").Value = "Tor Forsse"
$ws.Range("D2").Value = "Janne Tivenius"

$ws.Range("C3").Select()
